$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "3º Entrega" row: Aula time 1:30 -> 2:00 ---
$ws.Range("B15").Value = 0.083333333333333329

# --- Copy formatting of row 16 (last entry row) down into the two new rows ---
$ws.Range("A16:G16").Copy()
$ws.Range("A17:G18").PasteSpecial(-4122)

# Match the custom row height used by the rest of the table
$ws.Rows.Item(17).RowHeight = 20.1
$ws.Rows.Item(18).RowHeight = 20.1

# --- Row 17: new "5º Entrega" record ---
$ws.Range("A17").Value = "5º Entrega"
$ws.Range("B17").Value = 0.0625
$ws.Range("C17").Value = 0.041666666666666664

# --- Row 18: new "6º Entrega" record ---
$ws.Range("A18").Value = "6º Entrega"
$ws.Range("B18").Value = 0.083333333333333329
$ws.Range("C18").Value = 0.1875

# --- Restore the selected cell as left by the author ---
$ws.Range("N18").Select()
